# SPIKEY TREE: Added alternate graphic.
#
# Moves the existing "tree" shapes (clouds/stars/hexagon) and adds a new
# alternate "spikey" (star7) graphic cluster positioned where the green
# cloud cluster used to sit, built by duplicating the cloud shapes and
# switching their geometry to a 7-point star.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU_PER_PT = 12700.0

function EMU-ToPt([double]$v) {
    # The host stores Left/Top/Width/Height as points and truncates (rather
    # than rounds) when converting back to EMU internally, so a plain
    # v/12700.0 can land one EMU short because of float rounding (e.g.
    # 226244/12700.0 -> 17.814488188976377 -> *12700.0 ->
    # 226243.99999999997). Nudge by a tiny epsilon, far below a single EMU,
    # so the truncation recovers the exact target EMU value.
    return ($v / $EMU_PER_PT) + 0.00003
}

# --- 1. Shift the left-hand tree group (Cloud 5, Star 6, Star 4, Hexagon 3) up ---
$cloud5 = $s.Shapes.Item(1)   # "Cloud 5"
$star6  = $s.Shapes.Item(2)   # "Star: 7 Points 6"
$star4  = $s.Shapes.Item(3)   # "Star: 7 Points 4"
$hex3   = $s.Shapes.Item(4)   # "Hexagon 3"

$cloud5.Top = EMU-ToPt 226244
$star6.Top  = EMU-ToPt 226244
$star4.Top  = EMU-ToPt 457200
$hex3.Top   = EMU-ToPt 1175011

# --- 2. Shift the right-hand cloud cluster (Cloud 7, Cloud 9, Cloud 8) down/left ---
$cloud7 = $s.Shapes.Item(5)   # "Cloud 7"
$cloud9 = $s.Shapes.Item(6)   # "Cloud 9"
$cloud8 = $s.Shapes.Item(7)   # "Cloud 8"

$cloud7.Left = EMU-ToPt 742950
$cloud7.Top  = EMU-ToPt 3333750

$cloud9.Left = EMU-ToPt 1543050
$cloud9.Top  = EMU-ToPt 3814762

$cloud8.Left = EMU-ToPt 2063750
$cloud8.Top  = EMU-ToPt 4275398

# --- 3. Burn through the two low free shape ids (2,3) so the new shapes we ---
#        add below land on ids 11/12/13, matching the target document.
$dummy1 = $cloud7.Duplicate()
$dummy2 = $cloud7.Duplicate()
$dummy1.Delete()
$dummy2.Delete()

# --- 4. Add the new alternate "spikey tree" graphic: duplicates of the three ---
#        clouds, switched to 7-point star geometry and repositioned.
$msoShape7pointStar = 148

$star10 = $cloud7.Duplicate()
$star10.Name = "Star: 7 Points 10"
$star10.AutoShapeType = $msoShape7pointStar
$star10.Left = EMU-ToPt 6057902
$star10.Top  = EMU-ToPt 1304474

$star11 = $cloud9.Duplicate()
$star11.Name = "Star: 7 Points 11"
$star11.AutoShapeType = $msoShape7pointStar
$star11.Left = EMU-ToPt 6858002
$star11.Top  = EMU-ToPt 1785486

$star12 = $cloud8.Duplicate()
$star12.Name = "Star: 7 Points 12"
$star12.AutoShapeType = $msoShape7pointStar
$star12.Left = EMU-ToPt 7378702
$star12.Top  = EMU-ToPt 2246122
